# Generate Report for Handback
# - Update status text "Ready for handoff" -> "Handed back: in sync with en-US"
# - Add "Latest Target File" (F) / "Latest Handback File" (G) values + hyperlinks
#   for rows 2/3 on the zh-cn and de-de sheets (mirroring columns A and D).
# - Update "Latest Handback DateTime" (H) for rows 2/3 on zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$cornflower = 15570276  # RGB(0x64,0x95,0xED) packed as OLE BGR-ish value Excel expects

function Copy-LinkCell($ws, $srcAddr, $dstAddr) {
    # Find the hyperlink attached to the source cell (if any) and read its
    # target address + display text, then apply the same value/hyperlink/
    # formatting to the destination cell.
    $srcRange = $ws.Range($srcAddr)
    $dstRange = $ws.Range($dstAddr)

    $linkAddress = $null
    $linkText = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $srcRange.Address()) {
            $linkAddress = $h.Address
            $linkText = $h.TextToDisplay
            break
        }
    }

    if ($linkAddress -ne $null) {
        $dstRange.Value = $linkText
        $ws.Hyperlinks.Add($dstRange, $linkAddress, "", "", $linkText) | Out-Null
    } else {
        $dstRange.Value = $srcRange.Value()
    }

    $dstRange.Font.Color = $cornflower
    $dstRange.Font.Underline = 2
}

# ---- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet: new F/G columns + updated handback datetime ----
Copy-LinkCell $zhcn "A2" "F2"
Copy-LinkCell $zhcn "D2" "G2"
Copy-LinkCell $zhcn "A3" "F3"
Copy-LinkCell $zhcn "D3" "G3"

$zhcn.Range("H2").Value = "2016-03-11 09:35:49"
$zhcn.Range("H3").Value = "2016-03-11 09:35:49"

# ---- de-de sheet: new F/G columns + updated handback datetime ----
Copy-LinkCell $dede "A2" "F2"
Copy-LinkCell $dede "D2" "G2"
Copy-LinkCell $dede "A3" "F3"
Copy-LinkCell $dede "D3" "G3"

$dede.Range("H2").Value = "2016-03-11 09:35:54"
$dede.Range("H3").Value = "2016-03-11 09:35:54"

Write-Host "Handback report generated."
